# Problem Two: Potential Solutions
#
# Inserts the "Identify Potential Solutions" bullet and its four
# sub-bullets right after the "You need to select one matching pair of
# each color" bullet in the "Socks in the Dark" section, and before the
# trailing (bookmarked) empty bullet paragraph.

$d = $word.ActiveDocument

# --- locate the anchor paragraph -------------------------------------
# ("You need to select one matching pair of each color" - the last
# sub-bullet under "Break the Problem Apart")
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*You need to select one matching pair of each color*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -eq -1) {
    throw "Could not locate anchor paragraph."
}

# --- helpers -----------------------------------------------------------

# Inserts a brand new list paragraph right after the paragraph at
# $afterIndex, gives it the requested text + list level (0-based,
# matching w:ilvl) on numId 2 (the "Socks in the Dark" list), and
# returns its paragraph index.
function Add-ListParagraph($afterIndex, $ilvl, $text) {
    $anchorPara = $d.Paragraphs.Item($afterIndex)
    $r = $anchorPara.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = $text
    $newPara.Range.ListFormat.ListLevelNumber = $ilvl + 1
    return $newIndex
}

# Forces a run boundary at the start/end of the $occurrence-th occurrence
# of $substr inside paragraph $paraIndex by toggling (and then undoing) a
# character format on just that span - this leaves the visible text and
# formatting untouched but splits what would otherwise be one merged run
# into separate <w:r> runs, matching how Word naturally keeps such spans
# as distinct runs.
function Mark-RunBoundary($paraIndex, $substr, $occurrence) {
    $para = $d.Paragraphs.Item($paraIndex)
    $full = $para.Range.Text
    $pStart = $para.Range.Start
    $searchFrom = 0
    $off = -1
    for ($k = 0; $k -lt $occurrence; $k++) {
        $off = $full.IndexOf($substr, $searchFrom)
        $searchFrom = $off + 1
    }
    $s = $pStart + $off
    $e = $s + $substr.Length
    $rr = $d.Range($s, $e)
    $rr.Bold = 1
    $rr.Bold = 0
}

# --- insert the new bullets --------------------------------------------

$idx = $anchorIndex
$idx = Add-ListParagraph $idx 0 "Identify Potential Solutions"
$idx = Add-ListParagraph $idx 1 "You can select 4 socks and this would give a guaranteed match"
$idx = Add-ListParagraph $idx 1 "You can select 10 socks and this would only give you a guaranteed match with Black socks"
$idx16 = Add-ListParagraph $idx 1 "You can select 16 socks and this would only give you a guaranteed match with Black & Brown socks"
$idx18 = Add-ListParagraph $idx16 1 "You can select 18 socks and this would only give you a guaranteed match with Black & Brown & White socks"

# Re-create the original run boundaries. Must run in reverse document
# order (and right-to-left within a paragraph) to avoid perturbing runs
# that come later in the document.
Mark-RunBoundary $idx18 "& White " 1

Mark-RunBoundary $idx16 "Black & Brown" 1
Mark-RunBoundary $idx16 "guaranteed" 1
Mark-RunBoundary $idx16 "You can select 16" 1

Write-Output "Inserted Problem Two potential-solutions bullets."
